# The original workbook contains a header row (1), a real data row (2, John
# Glenn) and two more rows (3, 4) that were actually sample/template rows
# accidentally left in and committed ("discovered error in writing to
# original sample file"). This script clears those leftover sample rows
# back to blank cells (keeping the style-only placeholder cells + the
# number/text formats that were already on B3:B4, H3:J4), and drops the
# stray mailto hyperlink that only existed because of that sample data,
# matching the upstream fix-up commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the leftover sample-row values (row 3 "First Last" / row 4
# "KyleSant") while leaving any already-blank formatted cells alone.
$ws.Range("A3:G4").ClearContents()
$ws.Range("H3:J4").ClearContents()

# The H4 cell held a mailto: hyperlink tied to the sample row data; remove
# it along with the rest of that row's content.
$ws.Hyperlinks.Delete()

# Reflect the post-cleanup selection that was left active in the sheet.
$ws.Range("A3:X4").Select()
